# Auto-generated Excel COM-interop script to apply the diff
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Range("F2").Value = 902
$ws.Range("F3").Value = 14034
$ws.Range("F5").Value = 1066
$ws.Range("F6").Value = 822
$ws.Range("F8").Value = 643
$ws.Range("F10").Value = 30
$ws.Range("F11").Value = 69
$ws.Range("F12").Value = 790
$ws.Range("F13").Value = 2170
$ws.Range("F14").Value = 150
$ws.Range("F15").Value = 110
$ws.Range("F16").Value = 92
$ws.Range("F17").Value = 170
$ws.Range("F19").Value = 564
$ws.Range("F20").Value = 448
$ws.Range("F21").Value = 475
$ws.Range("F22").Value = 337
$ws.Range("F23").Value = 15
$ws.Range("F24").Value = 294
$ws.Range("F25").Value = 858
$ws.Range("F26").Value = 126
$ws.Range("F27").Value = 52
$ws.Range("F28").Value = 11
$ws.Range("F31").Value = 35
$ws.Range("F32").Value = 18
$ws.Range("C4").Value = "广州·珠三角 2024 COMIC WORLD次元世界动漫游戏嘉年华（取消）"

$ws = $wb.Worksheets.Item("演出")
$ws.Range("F6").Value = 95
$ws.Range("F7").Value = 181
$ws.Range("F8").Value = 1630
$ws.Range("F13").Value = 81
$ws.Range("F15").Value = 1724

$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F2").Value = 230
$ws.Range("F3").Value = 78
$ws.Range("F4").Value = 133

$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F2").Value = 230
$ws.Range("F3").Value = 902
$ws.Range("F4").Value = 14034
$ws.Range("F6").Value = 1066
$ws.Range("F7").Value = 822
$ws.Range("F9").Value = 643
$ws.Range("F11").Value = 30
$ws.Range("F12").Value = 69
$ws.Range("F13").Value = 790
$ws.Range("F16").Value = 2170
$ws.Range("F17").Value = 78
$ws.Range("F18").Value = 150
$ws.Range("F19").Value = 150
$ws.Range("F20").Value = 110
$ws.Range("F21").Value = 92
$ws.Range("F22").Value = 170
$ws.Range("F25").Value = 95
$ws.Range("F26").Value = 133
$ws.Range("F27").Value = 564
$ws.Range("F28").Value = 448
$ws.Range("F29").Value = 475
$ws.Range("F30").Value = 337
$ws.Range("F31").Value = 15
$ws.Range("F32").Value = 294
$ws.Range("F34").Value = 181
$ws.Range("F35").Value = 1630
$ws.Range("F40").Value = 126
$ws.Range("F41").Value = 52
$ws.Range("F42").Value = 11
$ws.Range("F44").Value = 81
$ws.Range("F47").Value = 35
$ws.Range("F48").Value = 18
$ws.Range("F49").Value = 1724
$ws.Range("C5").Value = "广州·珠三角 2024 COMIC WORLD次元世界动漫游戏嘉年华（取消）"
